$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '57.755.41'
$ws.Cells.Item(2, 5).Value = '  +1.89%  '
$ws.Cells.Item(3, 4).Value = '3.023.33'
$ws.Cells.Item(3, 5).Value = '  +0.69%  '
$ws.Cells.Item(4, 4).NumberFormat = "@"
$ws.Cells.Item(4, 4).Value = '1.00'
$ws.Cells.Item(4, 5).Value = '  -0.05%  '
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '510.29'
$ws.Cells.Item(5, 5).Value = '  +0.01%  '
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '140.27'
$ws.Cells.Item(6, 5).Value = '  +0.38%  '
$ws.Cells.Item(7, 5).Value = '  -0.05%  '
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = '0.436'
$ws.Cells.Item(8, 5).Value = '  +0.78%  '
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = '7.61'
$ws.Cells.Item(9, 5).Value = '  +0.53%  '
$ws.Cells.Item(10, 5).Value = '  +1.76%  '
$ws.Cells.Item(11, 5).Value = '  +3.61%  '
$ws.Cells.Item(12, 4).Value = '3.533.28'
$ws.Cells.Item(12, 5).Value = '  +0.47%  '
$ws.Cells.Item(13, 5).Value = '  +0.77%  '
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = '26.66'
$ws.Cells.Item(14, 5).Value = '  +3.67%  '
$ws.Cells.Item(15, 5).Value = '  +5.57%  '
$ws.Cells.Item(16, 5).Value = '  +5.65%  '
$ws.Cells.Item(17, 4).Value = '57.705.61'
$ws.Cells.Item(17, 5).Value = '  +1.68%  '
$ws.Cells.Item(18, 4).Value = '3.020.13'
$ws.Cells.Item(18, 5).Value = '  +0.67%  '
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '12.94'
$ws.Cells.Item(19, 5).Value = '  +3.64%  '
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '8.00'
$ws.Cells.Item(20, 5).Value = '  +1.85%  '
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '331.08'
$ws.Cells.Item(21, 5).Value = '  +0.48%  '
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '1.00'
$ws.Cells.Item(22, 5).Value = '  +0.07%  '
$ws.Cells.Item(23, 2).Value = 'LEO'
$ws.Cells.Item(23, 3).Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '5.74'
$ws.Cells.Item(23, 5).Value = '  -0.44%  '
$ws.Cells.Item(24, 2).Value = 'Polygon'
$ws.Cells.Item(24, 3).Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '0.501'
$ws.Cells.Item(24, 5).Value = '  +3.59%  '
$ws.Cells.Item(25, 2).Value = 'Litecoin'
$ws.Cells.Item(25, 3).Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = '64.83'
$ws.Cells.Item(25, 5).Value = '  +3.31%  '
$ws.Cells.Item(26, 2).Value = 'Kaspa'
$ws.Cells.Item(26, 3).Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '0.168'
$ws.Cells.Item(26, 5).Value = '  -2.67%  '
$ws.Cells.Item(27, 2).Value = 'Binance-PegBSC-USD'
$ws.Cells.Item(27, 3).Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '1.00'
$ws.Cells.Item(27, 5).Value = '  -0.35%  '
$ws.Cells.Item(28, 2).Value = 'PEPE'
$ws.Cells.Item(28, 3).Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Cells.Item(28, 4).Value = '0.0₃0924'
$ws.Cells.Item(28, 5).Value = '  +1.26%  '
$ws.Cells.Item(29, 2).Value = 'RenderToken'
$ws.Cells.Item(29, 3).Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = '6.85'
$ws.Cells.Item(29, 5).Value = '  +2.17%  '
$ws.Cells.Item(30, 2).Value = 'InternetComputer(DFINITY)'
$ws.Cells.Item(30, 3).Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = '7.39'
$ws.Cells.Item(30, 5).Value = '  +3.85%  '
$ws.Cells.Item(31, 2).Value = 'PancakeSwap'
$ws.Cells.Item(31, 3).Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = '1.81'
$ws.Cells.Item(31, 5).Value = '  +1.39%  '
$ws.Cells.Item(32, 2).Value = 'Fetch.AI'
$ws.Cells.Item(32, 3).Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '1.21'
$ws.Cells.Item(32, 5).Value = '  -4.49%  '
$ws.Cells.Item(33, 2).Value = 'EthereumClassic'
$ws.Cells.Item(33, 3).Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '20.73'
$ws.Cells.Item(33, 5).Value = '  +0.27%  '
$ws.Cells.Item(34, 2).Value = 'NEARProtocol'
$ws.Cells.Item(34, 3).Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '4.78'
$ws.Cells.Item(34, 5).Value = '  +4.37%  '
$ws.Cells.Item(35, 2).Value = 'Monero'
$ws.Cells.Item(35, 3).Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '155.06'
$ws.Cells.Item(35, 5).Value = '  -0.15%  '
$ws.Cells.Item(36, 2).Value = 'Aptos'
$ws.Cells.Item(36, 3).Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '5.91'
$ws.Cells.Item(36, 5).Value = '  +4.22%  '
$ws.Cells.Item(37, 2).Value = 'ImmutableX'
$ws.Cells.Item(37, 3).Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = '1.28'
$ws.Cells.Item(37, 5).Value = '  +0.60%  '
$ws.Cells.Item(38, 2).Value = 'EnergySwap'
$ws.Cells.Item(38, 3).Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = '24.52'
$ws.Cells.Item(38, 5).Value = '  +1.09%  '
$ws.Cells.Item(39, 2).Value = 'Hedera'
$ws.Cells.Item(39, 3).Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '0.0679'
$ws.Cells.Item(39, 5).Value = '  -0.78%  '
$ws.Cells.Item(40, 2).Value = 'RenzoRestakedETH'
$ws.Cells.Item(40, 3).Value = 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth'
$ws.Cells.Item(40, 4).Value = '3.053.25'
$ws.Cells.Item(40, 5).Value = '  +0.58%  '
$ws.Cells.Item(41, 2).Value = 'OKB'
$ws.Cells.Item(41, 3).Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '37.70'
$ws.Cells.Item(41, 5).Value = '  +2.04%  '
$ws.Cells.Item(42, 2).Value = 'Filecoin'
$ws.Cells.Item(42, 3).Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '3.88'
$ws.Cells.Item(42, 5).Value = '  +6.06%  '
$ws.Cells.Item(43, 2).Value = 'FirstDigitalUSD'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = '1.00'
$ws.Cells.Item(43, 5).Value = '  +0.03%  '
$ws.Cells.Item(44, 2).Value = 'Mantle'
$ws.Cells.Item(44, 3).Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '0.652'
$ws.Cells.Item(44, 5).Value = '  +0.50%  '
$ws.Cells.Item(45, 2).Value = 'Stacks'
$ws.Cells.Item(45, 3).Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '1.42'
$ws.Cells.Item(45, 5).Value = '  +0.25%  '
$ws.Cells.Item(46, 2).Value = 'Maker'
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Cells.Item(46, 4).Value = '2.234.72'
$ws.Cells.Item(46, 5).Value = '  -1.51%  '
$ws.Cells.Item(47, 2).Value = 'ONDO'
$ws.Cells.Item(47, 3).Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '0.990'
$ws.Cells.Item(47, 5).Value = '  -0.94%  '
$ws.Cells.Item(48, 2).Value = 'Cosmos'
$ws.Cells.Item(48, 3).Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '6.06'
$ws.Cells.Item(48, 5).Value = '  +4.32%  '
$ws.Cells.Item(49, 2).Value = 'VeChain'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '0.0239'
$ws.Cells.Item(49, 5).Value = '  -0.28%  '
$ws.Cells.Item(50, 2).Value = 'InjectiveProtocol'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '19.58'
$ws.Cells.Item(50, 5).Value = '  +1.22%  '
$ws.Cells.Item(51, 2).Value = 'dogwifhat'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '1.87'
$ws.Cells.Item(51, 5).Value = '  -6.34%  '
